$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet is protected; unprotect to allow writes, then re-protect at the end.
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A58).
$ws.Range("A58").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-19 for illustrative purposes only and are subject to change."

# Refreshed Weight (D) / Percent Change (E) values for holdings rows 2-55.
$updates = @{
    2 = @(0.01710541896968546, -0.005928853754940788)
    3 = @(0.04988388614021398, 0.0155119402640036)
    4 = @(0.01468680031770275, -0.001660123387548951)
    5 = @(0.00951816685536281, 0.001107726391581298)
    6 = @(0.01484666650983174, 0.01486455896778072)
    7 = @(0.01993022690785243, -0.005361305361305413)
    8 = @(0.004592495738603883, -0.00477099236641243)
    9 = @(0.006738498381997301, -0.005632811124801962)
    10 = @(0.01443303085938094, -0.01415598290598297)
    11 = @(0.008724766627319551, -0.0157099697885198)
    12 = @(0.01542385858603032, -0.01794411689310416)
    13 = @(0.003260090762613795, 0.02440749911567042)
    14 = @(0.00652018152522759, -0.009449694274596854)
    15 = @(0.01506356658047994, -0.01592134475103091)
    16 = @(0.01094813238265008, -0.01124052004333687)
    17 = @(0.02272656732953147, 0.01702766179540705)
    18 = @(0.00862545980475633, 0.005936146317985003)
    19 = @(0.01691919396599358, -0.002679628591013894)
    20 = @(0.01212044052691001, 0.001957266351329201)
    21 = @(0.007317008398774026, -0.008330556481172979)
    22 = @(0.01374242999643299, -0.01745435016111707)
    23 = @(0.01972709931624584, -0.01379932356257041)
    24 = @(0.01032902979606649, -0.01230342275670671)
    25 = @(0.02142681433136492, 0.0004797697105389798)
    26 = @(0.01108453923248808, -0.005680399500624267)
    27 = @(0.01954133559175631, 0.0058237661590963)
    28 = @(0.05559798225933367, -0.004480212395254424)
    29 = @(0.02059743432553831, -0.002898550724637627)
    30 = @(0.030451017058039, 0.01138281757195414)
    31 = @(0.01581775807632339, 0.009788053949903786)
    32 = @(0.01347449267689613, -0.01189370005575174)
    33 = @(0.02130793609121625, 0.007682030728123124)
    34 = @(0.03996009012341261, 0.002780333837949289)
    35 = @(0.01173889672953705, -0.003368137420006745)
    36 = @(0.009573388565587081, -0.003372843789149038)
    37 = @(0.01215305955621909, -0.003660024400162665)
    38 = @(0.007604187650751809, -0.002339789418952298)
    39 = @(0.01193204092126419, -0.01052901900359515)
    40 = @(0.01704360755657046, 0.003688524590164111)
    41 = @(0.01714798845035954, -0.02701518691588811)
    42 = @(0.03298970009669071, 0.009707865168539387)
    43 = @(0.01124209902860531, 0.002303151798640934)
    44 = @(0.02253829951794346, -0.06235837940723277)
    45 = @(0.01376374768528449, 0.0222126467736723)
    46 = @(0.007814036739306938, 0.06098440286896123)
    47 = @(0.01368101396549144, -0.02336570140454308)
    48 = @(0.01002590346309315, -0.007985803016859028)
    49 = @(0.01430749701931264, 0.008166028767634392)
    50 = @(0.008229352763994864, -0.006398039749523576)
    51 = @(0.0103518960650973, 0.01149008224479897)
    52 = @(0.008669413123037463, -0.0007525083612038852)
    53 = @(0.1430785286326874, 0.0001971220185295053)
    54 = @(0.04367292642313365, 0.001923816852635873)
    55 = @(1, -0.00121389444126252)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]
    $ws.Cells.Item($row, 5).Value = $vals[1]
}

# Restore sheet protection (structure/objects/scenarios locked, matching original settings).
$ws.Protect()
